# Update the "paises.xlsx" COVID-19 stats snapshot to the newer data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp shown in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 11:22"

# --- Three rankings swapped places because their totals changed; refresh the
#     country labels in column A to reflect the new sort order. ---
$ws.Range("A14").Value = "Belgica"
$ws.Range("A15").Value = "Brasil"

$ws.Range("A39").Value = "Indonesia"
$ws.Range("A40").Value = "Chequia"

$ws.Range("A81").Value = "Afganistan"
$ws.Range("A82").Value = "Cuba"
$ws.Range("A83").Value = "Ghana"

# --- Refresh the numeric statistics (Casos totales, Nuevos casos, Casos
#     activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#     rows whose figures moved. ---

# España (row 5)
$ws.Range("B5").Value = 204178
$ws.Range("C5").Value = 3968
$ws.Range("D5").Value = 82514
$ws.Range("E5").Value = 100382
$ws.Range("F5").Value = 7371
$ws.Range("G5").Value = 430
$ws.Range("H5").Value = 21282

# Belgica (row 14, new top-of-group figures)
$ws.Range("B14").Value = 40956
$ws.Range("C14").Value = 973
$ws.Range("D14").Value = 8895
$ws.Range("E14").Value = 26063
$ws.Range("F14").Value = 1071
$ws.Range("G14").Value = 170
$ws.Range("H14").Value = 5998

# Brasil (row 15, carried down from the old row 14)
$ws.Range("B15").Value = 40814
$ws.Range("C15").Value = 71
$ws.Range("D15").Value = 22991
$ws.Range("E15").Value = 15235
$ws.Range("F15").Value = 8318
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 2588

# Indonesia (row 39, new top-of-group figures)
$ws.Range("B39").Value = 7135
$ws.Range("C39").Value = 375
$ws.Range("D39").Value = 842
$ws.Range("E39").Value = 5677
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 26
$ws.Range("H39").Value = 616

# Chequia (row 40, carried down from the old row 39)
$ws.Range("B40").Value = 6914
$ws.Range("C40").Value = 14
$ws.Range("D40").Value = 1597
$ws.Range("E40").Value = 5121
$ws.Range("F40").Value = 75
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 196

# Australia (row 41)
$ws.Range("B41").Value = 6645
$ws.Range("C41").Value = 20
$ws.Range("E41").Value = 1889

# Moldavia (row 60)
$ws.Range("D60").Value = 505
$ws.Range("E60").Value = 1971
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 72

# Afganistan (row 81, new top-of-group figures)
$ws.Range("B81").Value = 1092
$ws.Range("C81").Value = 66
$ws.Range("D81").Value = 150
$ws.Range("E81").Value = 906
$ws.Range("F81").Value = 7

# Cuba (row 82, carried down from the old row 81)
$ws.Range("B82").Value = 1087
$ws.Range("D82").Value = 285
$ws.Range("E82").Value = 766
$ws.Range("F82").Value = 9
$ws.Range("H82").Value = 36

# Ghana (row 83, carried down from the old row 82)
$ws.Range("B83").Value = 1042
$ws.Range("D83").Value = 99
$ws.Range("E83").Value = 934
$ws.Range("F83").Value = 4
$ws.Range("H83").Value = 9

# Estado de Palestina (row 105)
$ws.Range("B105").Value = 461
$ws.Range("C105").Value = 12
$ws.Range("E105").Value = 386
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 4
